# Refresh the crypto price ("Price") and 1h volume change ("Volume(1h)")
# columns with the latest scraped values (GitHub Actions nightly update).
#
# Note: several "Price" values look like plain decimals (e.g. "310.45")
# but must stay TEXT, matching their original inline-string storage -
# assigning such a string straight to .Value lets Excel's type-inference
# auto-convert it to a number (dropping trailing zeros / switching to
# scientific notation for tiny values). Prefixing with a leading
# apostrophe forces Excel to keep it as literal text, exactly like a
# user typing '310.45 into the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.164.28'
$ws.Range("E2").Value = '  +0.91%  '
$ws.Range("D3").Value = '1.853.58'
$ws.Range("E3").Value = '  +0.54%  '
$ws.Range("E4").Value = '  +1.01%  '
$ws.Range("E5").Value = '  +1.04%  '
$ws.Range("D6").Value = "'" + '310.45'
$ws.Range("E6").Value = '  +0.46%  '
$ws.Range("E7").Value = '  +1.90%  '
$ws.Range("D8").Value = "'" + '0.3699'
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("D9").Value = "'" + '0.07273'
$ws.Range("E9").Value = '  +1.70%  '
$ws.Range("D10").Value = "'" + '0.9339'
$ws.Range("E10").Value = '  +0.80%  '
$ws.Range("D11").Value = "'" + '19.96'
$ws.Range("E11").Value = '  +2.01%  '
$ws.Range("D12").Value = "'" + '0.07812'
$ws.Range("E12").Value = '  +1.43%  '
$ws.Range("D13").Value = '1.838.66'
$ws.Range("E13").Value = '  +0.36%  '
$ws.Range("D14").Value = "'" + '5.399'
$ws.Range("E14").Value = '  +2.15%  '
$ws.Range("D15").Value = "'" + '6.507'
$ws.Range("E15").Value = '  +1.68%  '
$ws.Range("D16").Value = "'" + '89.67'
$ws.Range("E16").Value = '  +1.63%  '
$ws.Range("E17").Value = '  +1.02%  '
$ws.Range("D18").Value = "'" + '0.000008715'
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("E19").Value = '  +0.97%  '
$ws.Range("D20").Value = '27.198.51'
$ws.Range("E20").Value = '  +0.88%  '
$ws.Range("E21").Value = '  +1.39%  '
$ws.Range("D22").Value = "'" + '5.081'
$ws.Range("E22").Value = '  +1.12%  '
$ws.Range("D23").Value = "'" + '10.67'
$ws.Range("E23").Value = '  +0.49%  '
$ws.Range("D24").Value = "'" + '1.945'
$ws.Range("E24").Value = '  +0.58%  '
$ws.Range("D25").Value = "'" + '153.24'
$ws.Range("E25").Value = '  +0.66%  '
$ws.Range("D26").Value = "'" + '18.42'
$ws.Range("E26").Value = '  +1.05%  '
$ws.Range("D27").Value = "'" + '1.990'
$ws.Range("E27").Value = '  -1.59%  '
$ws.Range("D28").Value = "'" + '115.08'
$ws.Range("E28").Value = '  +0.76%  '
$ws.Range("D29").Value = "'" + '4.933'
$ws.Range("E29").Value = '  +1.12%  '
$ws.Range("D30").Value = "'" + '0.08876'
$ws.Range("E30").Value = '  +0.22%  '
$ws.Range("D31").Value = "'" + '3.305'
$ws.Range("E31").Value = '  +2.63%  '
$ws.Range("E32").Value = '  +0.44%  '
$ws.Range("D33").Value = "'" + '4.545'
$ws.Range("E33").Value = '  +1.67%  '
$ws.Range("D34").Value = "'" + '0.7381'
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("D35").Value = "'" + '2.692'
$ws.Range("E35").Value = '  -3.46%  '
$ws.Range("E36").Value = '  +3.09%  '
$ws.Range("D37").Value = "'" + '0.01993'
$ws.Range("E37").Value = '  +2.80%  '
$ws.Range("D38").Value = "'" + '0.05260'
$ws.Range("E38").Value = '  +1.42%  '
$ws.Range("D39").Value = "'" + '2.982'
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("E40").Value = '  +2.01%  '
$ws.Range("D41").Value = "'" + '7.046'
$ws.Range("E41").Value = '  +2.04%  '
$ws.Range("E42").Value = '  +0.55%  '
$ws.Range("D43").Value = "'" + '8.326'
$ws.Range("E43").Value = '  +2.48%  '
$ws.Range("D44").Value = "'" + '10.57'
$ws.Range("E44").Value = '  -0.25%  '
$ws.Range("E45").Value = '  +1.26%  '
$ws.Range("E46").Value = '  +1.09%  '
$ws.Range("D47").Value = "'" + '102.27'
$ws.Range("E47").Value = '  +1.82%  '
$ws.Range("D48").Value = "'" + '1.627'
$ws.Range("E48").Value = '  +1.23%  '
$ws.Range("D49").Value = "'" + '66.09'
$ws.Range("E49").Value = '  +0.91%  '
$ws.Range("D50").Value = "'" + '0.06061'
$ws.Range("E50").Value = '  +0.33%  '
$ws.Range("D51").Value = "'" + '0.8941'
$ws.Range("E51").Value = '  +0.34%  '
